$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "full_hybrid_rew"
$ws.Range("B1").Value = "full_hybrid_rew"
$ws.Range("C1").Value = "full_hybrid_rew"
$ws.Range("A2").Value = "'0.8888888888888888"
$ws.Range("B2").Value = "'0.8888888888888888"
$ws.Range("C2").Value = 0
$ws.Range("A3").Value = "'0.7608695652173914"
$ws.Range("B3").Value = "'0.8260869565217391"
$ws.Range("C3").Value = "'0.3333333333333333"
$ws.Range("A4").Value = "'0.8070175438596491"
$ws.Range("B4").Value = "'0.8070175438596491"
$ws.Range("C4").Value = "'0.16666666666666666"
$ws.Range("A5").Value = 0.65432098765432101
$ws.Range("B5").Value = 0.65432098765432101
$ws.Range("C5").Value = "'0.3333333333333333"
$ws.Range("A6").Value = "'0.4700854700854701"
$ws.Range("B6").Value = "'0.8461538461538461"
$ws.Range("C6").Value = "'0.6666666666666666"
$ws.Range("A7").Value = "'0.7311827956989247"
$ws.Range("B7").Value = "'0.7849462365591398"
$ws.Range("C7").Value = "'0.16666666666666666"
$ws.Range("A8").Value = "'0.24561403508771928"
$ws.Range("B8").Value = "'0.9210526315789473"
$ws.Range("C8").Value = "'0.6666666666666666"
$ws.Range("A9").Value = "'0.6222222222222222"
$ws.Range("B9").Value = "'0.6222222222222222"
$ws.Range("C9").Value = 0
$ws.Range("A10").Value = 0.875
$ws.Range("B10").Value = "'0.8942307692307693"
$ws.Range("C10").Value = "'0.16666666666666666"
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "'0.16666666666666666"
$ws.Range("A12").Value = "'0.6578947368421053"
$ws.Range("B12").Value = "'0.8245614035087719"
$ws.Range("C12").Value = "'0.3333333333333333"
$ws.Range("A13").Value = "'0.7843137254901961"
$ws.Range("B13").Value = "'0.8235294117647058"
$ws.Range("C13").Value = "'0.3333333333333333"
$ws.Range("A14").Value = "'0.5494505494505495"
$ws.Range("B14").Value = "'0.8021978021978022"
$ws.Range("C14").Value = 0.5
$ws.Range("A15").Value = "'0.8037383177570093"
$ws.Range("B15").Value = "'0.8037383177570093"
$ws.Range("C15").Value = "'0.16666666666666666"

$ws.Range("F7").Select()
